$d = $word.ActiveDocument

# Builds a standalone WordOpenXML "package" fragment containing exactly one
# list paragraph (numId 1 / ilvl 0, same indent as its siblings), optionally
# stamping the paragraph mark's run properties with <w:u w:val="none"/> (this
# is what Find/Replace + Font.Underline would leave behind on the pilcrow,
# but InsertXML lets us place it with surgical precision and nothing else).
function New-ParaPackageXml {
    param(
        [string]$Text,
        [bool]$Underlined
    )
    $pPrRpr = ""
    if ($Underlined) {
        $pPrRpr = "<w:rPr><w:u w:val=`"none`"/></w:rPr>"
    }
    $body = "<w:p><w:pPr><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"1`"/></w:numPr>" + `
            "<w:ind w:left=`"600`" w:hanging=`"360`"/>$pPrRpr</w:pPr>" + `
            "<w:r><w:rPr><w:rtl w:val=`"0`"/></w:rPr>" + `
            "<w:t xml:space=`"preserve`">$Text</w:t></w:r></w:p>"
    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
           '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
           '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
           '<w:body>' + $body + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    return $xml
}

# Locate a paragraph by its exact (trimmed) visible text.
function Find-ParagraphByText {
    param([string]$NeedleText)
    $count = $d.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $para = $d.Paragraphs.Item($i)
        $t = $para.Range.Text
        $t = $t.TrimEnd([char]13, [char]7)
        if ($t -eq $NeedleText) {
            return $para
        }
    }
    return $null
}

# 1) "SPOs: (SPOResource, Occurrence, Kind : Quads relative to container SPOs);"
#    -> "Kind: (Resource, Occurrence (Kind));"
#    (paragraph mark gains <w:u w:val="none"/>, run formatting untouched)
$p1 = Find-ParagraphByText "SPOs: (SPOResource, Occurrence, Kind : Quads relative to container SPOs);"
$p1.Range.InsertXML((New-ParaPackageXml "Kind: (Resource, Occurrence (Kind));" $true))

# 2) "Quad : (ContextResource, Occurrence, Kind : SPOs relative to container Quad);"
#    -> "CSPOs / Classes Kinds:"
#    (paragraph mark gains <w:u w:val="none"/>, run formatting untouched)
$p2 = Find-ParagraphByText "Quad : (ContextResource, Occurrence, Kind : SPOs relative to container Quad);"
$p2.Range.InsertXML((New-ParaPackageXml "CSPOs / Classes Kinds:" $true))

# 3) Insert new paragraph after it: "Subject: (SubjectKind, Occurrence (Quad));"
$p2 = Find-ParagraphByText "CSPOs / Classes Kinds:"
$r2 = $p2.Range
$r2.Collapse(0)
$r2.InsertParagraphAfter() | Out-Null
$p3 = $p2.Next()
$p3.Range.InsertXML((New-ParaPackageXml "Subject: (SubjectKind, Occurrence (Quad));" $false))

# 4) Insert another new paragraph after that: "Quad: (Context, Occurrence (Subject (Predicate (Object)));"
$p3 = Find-ParagraphByText "Subject: (SubjectKind, Occurrence (Quad));"
$r3 = $p3.Range
$r3.Collapse(0)
$r3.InsertParagraphAfter() | Out-Null
$p4 = $p3.Next()
$p4.Range.InsertXML((New-ParaPackageXml "Quad: (Context, Occurrence (Subject (Predicate (Object)));" $false))

Write-Host "Done. Paragraph count: $($d.Paragraphs.Count)"
